$p = $ppt.ActivePresentation

$m1 = $p.SlideMaster
$cs1 = $m1.ColorScheme

# Target "Office Theme" palette (12-slot ColorScheme order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink), expressed as
# OLE_COLOR (0x00BBGGRR) longs.
$officeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

for ($i = 1; $i -le 12; $i++) {
    $c1 = $cs1.Colors($i)
    $c1.RGB = $officeColors[$i - 1]
}
